# Dummy Contacts workbook update
# - Change the email for the first "Sam/Saun" contact row (row 2)
# - Rename the last contact row (row 10) to a "Row"/"ShouldFail" test row
# - Add a new "Statuscode" column (K) marking each row Active/Inactive
# - Update the active cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new Statuscode column
$ws.Range("K1").Value = "Statuscode"
$ws.Range("K2").Value = "Active"
$ws.Range("K3").Value = "Active"
$ws.Range("K4").Value = "Active"
$ws.Range("K5").Value = "Active"
$ws.Range("K6").Value = "Active"
$ws.Range("K7").Value = "Active"
$ws.Range("K8").Value = "Active"
$ws.Range("K9").Value = "Active"
$ws.Range("K10").Value = "Inactive"

# Turn row 10 into a dedicated "should fail" validation test row
$ws.Range("A10").Value = "Row"
$ws.Range("B10").Value = "ShouldFail"

# Update email address on row 2
$ws.Range("C2").Value = "Sammmy.sean@samsam.com"

# Match the saved selection state
$ws.Range("E2").Select() | Out-Null
